# PBARC Rotating Olfactometers SOP - "info about circuit breaker" edit
$d = $word.ActiveDocument
$NL = [char]11   # vertical-tab -> manual line break (<w:br/>) when typed/found by Word

# ---------------------------------------------------------------------------
# 1) Add a sentence about the circuit breaker right after the paragraph that
#    ends "...is going to be unused for a while."
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("is going to be unused for a while.") | Out-Null
$ip1 = $d.Range($r1.End, $r1.End)
$ip1.InsertAfter("  Directly above the on/off switch is a push-to-reset circuit breaker.")

# ---------------------------------------------------------------------------
# 2) Insert a new bullet line "Has the circuit breaker tripped? ..." right
#    after "Is the power unit switched on?" (before the "Are the wires..."
#    line). Use the paragraph-split/rejoin trick so the new sentence lands in
#    its own run instead of being folded into the previous run.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Is the power unit switched on?") | Out-Null
$afterQuestion = $r2.End
$ip2 = $d.Range($afterQuestion, $afterQuestion)
$ip2.InsertParagraphAfter()

$r3 = $d.Content
$r3.Find.Execute("Is the power unit switched on?" + [char]13) | Out-Null
$afterMark = $r3.End
$ip3 = $d.Range($afterMark, $afterMark)
$ip3.InsertAfter($NL + "Has the circuit breaker tripped?  Push the button above the on/off switch.")

# remove the temporary paragraph mark that separated the two sentences so they
# end up back in a single paragraph but in distinct runs
$markRange = $d.Range($afterQuestion, $afterMark)
$markRange.Delete()

# ---------------------------------------------------------------------------
# 3) Reword "The switch on the power unit isn't very good and might have
#    failed." -> "The actual on/off switch on the power unit isn't very
#    robust and might have failed."
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("good and might have failed.", $true, $false, $false, $false, $false, $true, 1, $false, "robust and might have failed.", 2) | Out-Null

$r5 = $d.Content
$r5.Find.Execute("The switch on the power unit") | Out-Null
$ip5 = $d.Range($r5.Start + 3, $r5.Start + 3)
$ip5.InsertAfter(" actual on/off")
